$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A
$ws.Columns("A").Insert()

# Set header
$ws.Range("A1").Value = "network_code"
$ws.Range("A1").Font.Bold = $true

# Fill A2:A16 with "KJ"
$ws.Range("A2:A16").Value = "KJ"

# Autofit the new column to match the content width
$ws.Columns("A").AutoFit()

# Match the final selection left by the author's session
$ws.Range("P19").Select() | Out-Null
